# Update the "dSF" column (column F) for several rows with re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -2
    11 = -6
    12 = -2
    23 = -4
    25 = -5
    26 = -8
    35 = -8
    41 = -4
    45 = 0
    46 = -1
    50 = 6
    53 = -2
    55 = 6
    59 = -2
    61 = 5
    70 = -1
    71 = -2
    76 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
